$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "AVETP651"
$ws.Range("B2").Value = 23101918
$ws.Range("C2").Value = "vesdxjm16"
$ws.Range("D2").Value = "D39fQg%#"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "qjsaBDOH"
$ws.Range("G2").Value = "vTEE"
$ws.Range("H2").Value = "Candidate"

# Delete rows 3 and 4 (the two extra candidate rows)
$ws.Range("A3:H4").EntireRow.Delete()
